# Apply the edit described by the diff:
#  - Insert a new row at position 4 on the first worksheet (strategy_id-0),
#    shifting the existing rows 4-11 down to rows 5-12.
#  - Populate the new row 4 with the "climate_change_factor_gnrl_hydropower_availability"
#    variable and its sampling range values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Insert a new row before the current row 4 ("elasticity_gnrl_rate_occupancy_to_gdppc"),
# pushing it (and everything below) down by one row.
$ws.Rows.Item(4).Insert()

# Fill in the metadata columns for the newly inserted row 4.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

# max_35 / min_35 style columns.
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

# Yearly sampling values (columns J through AS) are all set to 1.
$ws.Range("J4:AS4").Value = 1

